# NYPD 112th Precinct CompStat weekly report refresh:
# new week's crime data collected -> bump volume/week numbers, update the
# weekly crime-complaint stats table (rows 15-28), and widen column E to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header: bump the bulletin "Number" and the reporting week date range
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/16/2024  Through  9/22/2024"

# ---------------------------------------------------------------------------
# Column E ("Week to Date % Chg") needs to be a bit wider this week
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 7.433768

# ---------------------------------------------------------------------------
# Helper: convert a "N/A"/"***.*" placeholder-text cell into a real number,
# keeping the normal right-aligned number style used throughout the table.
# ---------------------------------------------------------------------------
function Set-NumberCell($addr, $value, $format) {
    $rng = $ws.Range($addr)
    $rng.Value = $value
    $rng.NumberFormat = $format
}

$fmtCount = "#,##0"
$fmtPct = "#,##0.0;""-""#,##0.0"

# ---------------------------------------------------------------------------
# Helper: convert a real number cell into a "N/A"/"***.*" placeholder-text
# cell, copying the donor cell's format (General, same font) then writing
# the literal text as a true string (not a numeric-looking value).
# ---------------------------------------------------------------------------
function Set-PlaceholderCell($addr, $text, $donor) {
    $ws.Range($donor).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Formula = "=""" + $text + """"
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 8
$ws.Range("K15").Value = 166.666666666667
$ws.Range("L15").Value = -11.111111111111
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 0

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -63.636363636363
$ws.Range("I16").Value = 60
$ws.Range("J16").Value = 54
$ws.Range("K16").Value = 11.111111111111
$ws.Range("L16").Value = -10.447761194029
$ws.Range("M16").Value = -24.050632911392
$ws.Range("N16").Value = -87.179487179487

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 75
$ws.Range("J17").Value = 73
$ws.Range("K17").Value = 2.739726027397
$ws.Range("L17").Value = 13.636363636363
$ws.Range("M17").Value = 82.926829268292
$ws.Range("N17").Value = -14.772727272727

# ---------------------------------------------------------------------------
# Row 18 - Burglary (D18/E18 flip from "N/A"/"***.*" to real numbers)
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 1
Set-NumberCell "D18" 2 $fmtCount
Set-NumberCell "E18" -50 $fmtPct
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 70
$ws.Range("J18").Value = 85
$ws.Range("K18").Value = -17.647058823529
$ws.Range("L18").Value = -14.634146341463
$ws.Range("M18").Value = -23.076923076923
$ws.Range("N18").Value = -93.006993006993

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 0
$ws.Range("G19").Value = 33
$ws.Range("H19").Value = -3.030303030303
$ws.Range("I19").Value = 320
$ws.Range("J19").Value = 358
$ws.Range("K19").Value = -10.614525139664
$ws.Range("L19").Value = -14.893617021276
$ws.Range("M19").Value = 16.788321167883
$ws.Range("N19").Value = -56.284153005464

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -80
$ws.Range("I20").Value = 125
$ws.Range("J20").Value = 110
$ws.Range("K20").Value = 13.636363636363
$ws.Range("L20").Value = 73.611111111111
$ws.Range("M20").Value = 71.232876712328
$ws.Range("N20").Value = -95.015948963317

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -28.571428571428
$ws.Range("F21").Value = 64
$ws.Range("G21").Value = 75
$ws.Range("H21").Value = -14.666666666666
$ws.Range("I21").Value = 658
$ws.Range("J21").Value = 683
$ws.Range("K21").Value = -3.660322108345
$ws.Range("L21").Value = -2.373887240356
$ws.Range("M21").Value = 17.290552584670
$ws.Range("N21").Value = -86.314475873544

# ---------------------------------------------------------------------------
# Row 22 - Transit (D22/E22 flip from "N/A"/"***.*" to real numbers)
# ---------------------------------------------------------------------------
$ws.Range("C22").Value = 2
Set-NumberCell "D22" 1 $fmtCount
Set-NumberCell "E22" 100 $fmtPct
$ws.Range("I22").Value = 28
$ws.Range("J22").Value = 23
$ws.Range("K22").Value = 21.739130434782
$ws.Range("L22").Value = 7.692307692307
$ws.Range("M22").Value = 75

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 41
$ws.Range("E24").Value = -19.512195121951
$ws.Range("F24").Value = 134
$ws.Range("G24").Value = 134
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 1246
$ws.Range("J24").Value = 1125
$ws.Range("K24").Value = 10.755555555555
$ws.Range("L24").Value = -5.534495830174
$ws.Range("M24").Value = 74.754558204768

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 26
$ws.Range("D25").Value = 32
$ws.Range("E25").Value = -18.75
$ws.Range("F25").Value = 101
$ws.Range("G25").Value = 106
$ws.Range("H25").Value = -4.716981132075
$ws.Range("I25").Value = 905
$ws.Range("J25").Value = 801
$ws.Range("K25").Value = 12.983770287141
$ws.Range("L25").Value = -3.620873269435

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -30
$ws.Range("F26").Value = 21
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = -25
$ws.Range("I26").Value = 206
$ws.Range("J26").Value = 183
$ws.Range("K26").Value = 12.568306010929
$ws.Range("L26").Value = 45.070422535211
$ws.Range("M26").Value = 34.640522875817

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape* (D27/E27 flip from real numbers to "N/A"/"***.*")
# ---------------------------------------------------------------------------
$ws.Range("C27").Value = 1
Set-PlaceholderCell "D27" "0" "D23"
Set-PlaceholderCell "E27" "***.*" "D23"
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 16
$ws.Range("K27").Value = 77.777777777777
$ws.Range("L27").Value = 14.285714285714

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes (C28 flips from "N/A" to a real number)
# ---------------------------------------------------------------------------
Set-NumberCell "C28" 1 $fmtCount
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 24
$ws.Range("K28").Value = 4.347826086956
$ws.Range("L28").Value = -25
